$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.420.86'
$ws.Range("E2").Value = '  +1.26%  '

$ws.Range("D3").Value = '1.677.31'
$ws.Range("E3").Value = '  +2.42%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.93'
$ws.Range("E5").Value = '  +1.32%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5313'
$ws.Range("E6").Value = '  +0.98%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2698'
$ws.Range("E8").Value = '  +3.76%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06403'
$ws.Range("E9").Value = '  +1.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.74'
$ws.Range("E10").Value = '  +4.85%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07817'
$ws.Range("E11").Value = '  +2.10%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.682.89'
$ws.Range("E12").Value = '  +2.54%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.513'
$ws.Range("E13").Value = '  +1.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5566'
$ws.Range("E14").Value = '  +1.07%  '

$ws.Range("D15").Value = '0.0₅8341'
$ws.Range("E15").Value = '  +1.81%  '

$ws.Range("E16").Value = '  +0.77%  '

$ws.Range("D17").Value = '26.483.67'
$ws.Range("E17").Value = '  +1.55%  '

$ws.Range("E18").Value = '  -0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.732'
$ws.Range("E19").Value = '  +0.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.58'
$ws.Range("E20").Value = '  +2.81%  '

$ws.Range("E21").Value = '  +1.35%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.340'
$ws.Range("E22").Value = '  +2.78%  '

$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '142.23'
$ws.Range("E24").Value = '  -2.56%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1289'
$ws.Range("E25").Value = '  +5.79%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.406'
$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("E27").Value = '  +2.51%  '

$ws.Range("E28").Value = '  +1.90%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06339'
$ws.Range("E29").Value = '  +5.40%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.273'
$ws.Range("E30").Value = '  +1.34%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.614'
$ws.Range("E31").Value = '  +4.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.448'
$ws.Range("E32").Value = '  +1.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.677'
$ws.Range("E33").Value = '  +2.21%  '

$ws.Range("E34").Value = '  +2.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6154'
$ws.Range("E35").Value = '  +7.15%  '

$ws.Range("E36").Value = '  +1.33%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.783'
$ws.Range("E37").Value = '  +0.64%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.163'
$ws.Range("E38").Value = '  +7.66%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01632'
$ws.Range("E39").Value = '  +0.63%  '

$ws.Range("D40").Value = '1.084.87'
$ws.Range("E40").Value = '  +4.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8647'
$ws.Range("E41").Value = '  +1.14%  '

$ws.Range("E42").Value = '  -0.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.37'

$ws.Range("D44").Value = '1.823.03'
$ws.Range("E44").Value = '  +1.95%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.31'
$ws.Range("E45").Value = '  +3.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.178'
$ws.Range("E46").Value = '  +1.39%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.001'
$ws.Range("E47").Value = '  -0.08%  '

$ws.Range("E48").Value = '  -3.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05207'
$ws.Range("E49").Value = '  +0.71%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.472'
$ws.Range("E50").Value = '  +5.75%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.031'
$ws.Range("E51").Value = '  +1.78%  '
